$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (F1, G1, X1, Y1, Z1) ---
$ws.Range("F1").Value = "FT_Goals_H"
$ws.Range("G1").Value = "FT_Goals_A"
$ws.Range("X1").Value = "FT_Odds_H"
$ws.Range("Y1").Value = "FT_Odds_D"
$ws.Range("Z1").Value = "FT_Odds_A"

# --- Append 9 new match rows (145-153) ---
$dataRows = @(
    @('T1', '03/01/2023', '14:00', 'Umraniyespor', 'Hatayspor', 2, 2, 'D', 1, 1, 'D', 11, 17, 3, 5, 10, 14, 5, 2, 3, 2, 1, 0, 2.2, 3.3, 3.1, 1.83, 2.02),
    @('T1', '03/01/2023', '17:00', 'Antalyaspor', 'Fenerbahce', 1, 2, 'A', 1, 0, 'H', 8, 18, 2, 4, 21, 17, 3, 8, 4, 2, 0, 0, 4.2, 3.75, 1.75, 1.7, 2.1),
    @('T1', '04/01/2023', '14:00', 'Kasimpasa', 'Alanyaspor', 4, 1, 'H', 2, 0, 'H', 11, 14, 6, 2, 14, 13, 1, 7, 2, 0, 0, 0, 2.25, 3.3, 2.9, 1.62, 2.25),
    @('T1', '04/01/2023', '14:00', 'Kayserispor', 'Karagumruk', 2, 4, 'A', 0, 2, 'A', 19, 11, 5, 7, 12, 16, 6, 2, 1, 3, 0, 0, 1.91, 3.5, 3.6, 1.62, 2.25),
    @('T1', '04/01/2023', '17:00', 'Galatasaray', 'Ankaragucu', 2, 1, 'H', 2, 1, 'H', 12, 9, 5, 2, 17, 15, 6, 4, 2, 3, 0, 0, 1.44, 4.2, 6.5, 1.57, 2.35),
    @('T1', '04/01/2023', '17:00', 'Konyaspor', 'Sivasspor', 2, 2, 'D', 1, 1, 'D', 10, 10, 5, 3, 13, 13, 4, 3, 4, 1, 0, 0, 1.85, 3.4, 4, 2, 1.85),
    @('T1', '05/01/2023', '14:00', 'Gaziantep', 'Buyuksehyr', 1, 1, 'D', 0, 1, 'A', 14, 16, 3, 6, 8, 6, 5, 5, 2, 0, 0, 0, 3.25, 3.3, 2.1, 1.95, 1.9),
    @('T1', '05/01/2023', '17:00', 'Ad. Demirspor', 'Istanbulspor', 6, 0, 'H', 1, 0, 'H', 17, 8, 11, 3, 12, 11, 4, 5, 1, 3, 0, 0, 1.4, 4.33, 6.5, 1.53, 2.4),
    @('T1', '05/01/2023', '17:00', 'Trabzonspor', 'Giresunspor', 3, 0, 'H', 1, 0, 'H', 17, 5, 8, 1, 11, 9, 10, 1, 2, 2, 0, 0, 1.45, 4.2, 6.5, 1.75, 2.05),
)

$startRow = 145
$numRows = $dataRows.Count

# Format the Date/Time columns (B,C) as Text first so strings like
# "03/01/2023" / "14:00" are not auto-converted into date/time serials.
$dateTimeRange = $ws.Range("B" + $startRow + ":C" + ($startRow + $numRows - 1))
$dateTimeRange.NumberFormat = "@"

for ($i = 0; $i -lt $numRows; $i++) {
    $r = $startRow + $i
    $rowData = $dataRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

# Restore the default (Normal) style on the Date/Time cells now that the
# text values are safely in place, so no stray number-format style lingers.
$dateTimeRange.Style = "Normal"
